$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "userId" header in column E to "_id"
$ws.Range("E1").Value = "_id"

# Populate the new _id values for each user row (previously column E was empty)
$ws.Range("E2").Value = "aa847edee5847831acb269a4"
$ws.Range("E3").Value = "aa847edee5847831acb269a5"
$ws.Range("E4").Value = "aa847edee5847831acb269a6"
$ws.Range("E5").Value = "aa847edee5847831acb269a7"
$ws.Range("E6").Value = "aa847edee5847831acb269a8"
$ws.Range("E7").Value = "aa847edee5847831acb269a9"
$ws.Range("E8").Value = "aa847edee5847831acb269aa"

# Widen the new column to fit its content (OOXML stores this as width="29")
$ws.Columns.Item(5).ColumnWidth = 28.14

# Match the author's final selection/active cell
$ws.Range("E8").Select()

$ws.PageSetup.Orientation = 1
